$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell forcing text storage even when the
# string looks like a plain number (e.g. "22.05"), so Excel does not
# silently coerce it to a numeric cell. Cells whose text is unambiguous
# (contains thousand-separator dots, subscript digits, etc.) are written
# directly since Excel already stores those as text.
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.041.34"
$ws.Range("E2").Value = "  +0.27%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.562.80"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.38%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "208.38"
$ws.Range("E5").Value = "  +0.53%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.27%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.41%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "22.05"
$ws.Range("E8").Value = "  -0.43%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.83%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0598"
$ws.Range("E10").Value = "  +1.59%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.23%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.784.12"
$ws.Range("E12").Value = "  +0.30%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.565.89"
$ws.Range("E13").Value = "  +0.58%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.14%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.06%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "27.059.06"
$ws.Range("E16").Value = "  +0.33%  "

# Row 17 - Litecoin
$ws.Range("E17").Value = "  +0.18%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0707"
$ws.Range("E18").Value = "  +1.40%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "216.72"
$ws.Range("E19").Value = "  -0.87%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "7.38"
$ws.Range("E20").Value = "  +0.72%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.43%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.70%  "

# Row 23 - Avalanche
Set-TextValue $ws.Range("D23") "9.22"
$ws.Range("E23").Value = "  -0.35%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.18%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "153.04"
$ws.Range("E25").Value = "  -0.76%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "6.60"
$ws.Range("E26").Value = "  -0.77%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "15.11"
$ws.Range("E27").Value = "  +0.81%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.32%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.36%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.92%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.77%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.25%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "3.16"
$ws.Range("E33").Value = "  +2.45%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.428.32"
$ws.Range("E34").Value = "  +0.18%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +1.50%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +7.92%  "

# Row 37 - HuobiToken
Set-TextValue $ws.Range("D37") "2.34"
$ws.Range("E37").Value = "  +2.47%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.20%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +2.13%  "

# Row 40 - FraxShare
Set-TextValue $ws.Range("D40") "5.89"
$ws.Range("E40").Value = "  +2.39%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  -0.47%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.45%  "

# Row 43 - now WEMIXToken (swapped with row 44)
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +1.33%  "

# Row 44 - now MXToken (swapped with row 43)
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D44") "2.32"
$ws.Range("E44").Value = "  -0.61%  "

# Row 45 - Aave
Set-TextValue $ws.Range("D45") "64.69"
$ws.Range("E45").Value = "  +0.13%  "

# Row 46 - RenderToken
Set-TextValue $ws.Range("D46") "1.75"
$ws.Range("E46").Value = "  -1.11%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.700.56"
$ws.Range("E47").Value = "  +0.45%  "

# Row 48 - Quant
Set-TextValue $ws.Range("D48") "86.91"
$ws.Range("E48").Value = "  -1.24%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +5.24%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -0.34%  "

# Row 51 - Algorand
Set-TextValue $ws.Range("D51") "0.0959"
$ws.Range("E51").Value = "  +0.22%  "
